$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header stays the same (SENTENCES / NAME) ---

# --- Rows 2-9: replace the old placeholder "Test#" rows with the
#     real comparison-sentence text, and renumber the P1_W1/P1_W2 labels ---
$ws.Range("A2").Value = "We picked grapes for wine"
$ws.Range("B2").Value = "P1_W1_S1"

$ws.Range("A3").Value = "The ballet is about to begin."
$ws.Range("B3").Value = "P1_W1_S2"

$ws.Range("A4").Value = "You're used to being on the field."
$ws.Range("B4").Value = "P1_W1_S3"

$ws.Range("A5").Value = "Enjoy the fair weather while in the tropics."
$ws.Range("B5").Value = "P1_W1_S4"

$ws.Range("A6").Value = "I think I'm getting better."
$ws.Range("B6").Value = "P1_W2_S1"

$ws.Range("A7").Value = "You want him to do well"
$ws.Range("B7").Value = "P1_W2_S2"

# --- New rows 8-9 (week 2 sentences that did not exist before) ---
$ws.Range("A8").Value = "Big muscles are not necessarily strong ones"
$ws.Range("B8").Value = "P1_W2_S3"

$ws.Range("A9").Value = "he is capable and willing to make decisions."
$ws.Range("B9").Value = "P1_W2_S4"

# --- Formatting: the whole sentence column gets an explicit black Calibri font ---
$ws.Range("A2:A9").Font.Name = "Calibri"
$ws.Range("A2:A9").Font.Color = 0

# --- The newly-added "week 2" labels (B6:B9) get the same explicit font ---
$ws.Range("B6:B9").Font.Name = "Calibri"
$ws.Range("B6:B9").Font.Color = 0

# --- Column A is now much wider to fit the long sentences ---
$ws.Columns.Item(1).ColumnWidth = 33.4

# --- Selection moved to A11 (matches the saved view state in the workbook) ---
$null = $ws.Range("A11").Select()
